# miniWECC_split.xlsx update: "update of miniWECC split plus weekly"
# Adds a new "Area" column (F) to the split-table starting at row 49, a
# "Nukes"/"res"/"*" annotation column (I), and a second mini Name/Area/Bus
# summary block in columns N:S for several rows. Also updates a couple of
# "*" flags (G53, G75) to " -- no gov" and G76 to " --" plus moves the old
# F76 "Nukes" note to I76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-seed the two brand-new shared strings in the same order they were
# first introduced upstream (" --" before " -- no gov") so the exported
# sharedStrings.xml allocates the same indices (138, 139) as the target.
$ws.Range("G76").Value = " --"
$ws.Range("G53").Value = " -- no gov"

# --- Header row additions ---
$ws.Range("F49").Value = "Area"

# --- Row 50 ---
$ws.Range("F50").Value = 1
$ws.Range("N50").Value = 1
$ws.Range("O50").Value = "BCH-G1      "
$ws.Range("P50").Value = 1
$ws.Range("R50").Value = 1
$ws.Range("S50").Value = "*"

# --- Row 51 ---
$ws.Range("F51").Value = 1

# --- Row 52 ---
$ws.Range("F52").Value = 1
$ws.Range("N52").Value = 17
$ws.Range("O52").Value = "WA-GEN      "
$ws.Range("P52").Value = 7
$ws.Range("R52").Value = 1
$ws.Range("S52").Value = "*"

# --- Row 53 ---
$ws.Range("F53").Value = 1
$ws.Range("N53").Value = 23
$ws.Range("O53").Value = "ORE-G23     "
$ws.Range("P53").Value = 10
$ws.Range("R53").Value = 1
$ws.Range("S53").Value = "*"

# --- Row 54 ---
$ws.Range("F54").Value = 1
$ws.Range("N54").Value = 118
$ws.Range("O54").Value = "ALB-GEN    "
$ws.Range("P54").Value = 34
$ws.Range("R54").Value = 1
$ws.Range("S54").Value = "*"

# --- Row 55 ---
$ws.Range("F55").Value = 1

# --- Row 56 ---
$ws.Range("F56").Value = 1
$ws.Range("I56").Value = "res"

# --- Row 57 ---
$ws.Range("F57").Value = 1
$ws.Range("I57").Value = "res"

# --- Row 58 ---
$ws.Range("F58").Value = 1
$ws.Range("I58").Value = "res"

# --- Row 59 ---
$ws.Range("F59").Value = 1
$ws.Range("I59").Value = "res"

# --- Row 60 ---
$ws.Range("F60").Value = 1
$ws.Range("I60").Value = "res"

# --- Row 61 ---
$ws.Range("F61").Value = 1
$ws.Range("N61").Value = 68
$ws.Range("O61").Value = "COLO-GEN    "
$ws.Range("P61").Value = 29
$ws.Range("R61").Value = 2
$ws.Range("S61").Value = "*"

# --- Row 62 ---
$ws.Range("F62").Value = 2
$ws.Range("N62").Value = 71
$ws.Range("O62").Value = "COLO-GEN    "
$ws.Range("P62").Value = 30
$ws.Range("R62").Value = 2
$ws.Range("S62").Value = "*"

# --- Row 63 ---
$ws.Range("F63").Value = 2

# --- Row 64 ---
$ws.Range("F64").Value = 1

# --- Row 65 ---
$ws.Range("F65").Value = 3

# --- Row 66 ---
$ws.Range("F66").Value = 3

# --- Row 67 ---
$ws.Range("F67").Value = 3

# --- Row 68 ---
$ws.Range("F68").Value = 3

# --- Row 69 ---
$ws.Range("F69").Value = 3

# --- Row 70 ---
$ws.Range("F70").Value = 3
$ws.Range("N70").Value = 41
$ws.Range("O70").Value = "SFO-GEN     "
$ws.Range("P70").Value = 18
$ws.Range("R70").Value = 3
$ws.Range("S70").Value = "*"

# --- Row 71 ---
$ws.Range("F71").Value = 3
$ws.Range("N71").Value = 48
$ws.Range("O71").Value = "SC-G1       "
$ws.Range("P71").Value = 21
$ws.Range("R71").Value = 3
$ws.Range("S71").Value = "*"

# --- Row 72 ---
$ws.Range("F72").Value = 3
$ws.Range("N72").Value = 59
$ws.Range("O72").Value = "NEV-G2      "
$ws.Range("P72").Value = 25
$ws.Range("R72").Value = 3
$ws.Range("S72").Value = "*"

# --- Row 73 ---
$ws.Range("F73").Value = 3

# --- Row 74 ---
$ws.Range("F74").Value = 3
$ws.Range("N74").Value = 65
$ws.Range("O74").Value = "AZ-65       "
$ws.Range("P74").Value = 28
$ws.Range("R74").Value = 3
$ws.Range("S74").Value = "*"

# --- Row 75 ---
$ws.Range("F75").Value = 3
$ws.Range("G75").Value = " -- no gov"

# --- Row 76 ---
$ws.Range("F76").Value = 3
$ws.Range("I76").Value = "Nukes"

# --- Row 77 ---
$ws.Range("F77").Value = 3

# --- Row 78 ---
$ws.Range("F78").Value = 2

# --- Row 79 ---
$ws.Range("F79").Value = 2

# --- Row 80 ---
$ws.Range("F80").Value = 2

# --- Row 81 ---
$ws.Range("F81").Value = 1

# --- Row 82 ---
$ws.Range("F82").Value = 2

# --- Row 83 ---
$ws.Range("F83").Value = 1

# --- Selection / view state ---
$ws.Range("G54").Select()
